$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HMI Internal")

$ws.Range("A146:A147").EntireRow.Delete()
$ws.Range("A140:A141").EntireRow.Delete()
$ws.Range("A134:A135").EntireRow.Delete()
$ws.Range("A128:A129").EntireRow.Delete()
$ws.Range("A122:A123").EntireRow.Delete()
$ws.Range("A116:A117").EntireRow.Delete()
$ws.Range("A110:A111").EntireRow.Delete()
$ws.Range("A104:A105").EntireRow.Delete()
$ws.Range("A98:A99").EntireRow.Delete()
$ws.Range("A92:A93").EntireRow.Delete()
$ws.Range("A86:A87").EntireRow.Delete()
$ws.Range("A80:A81").EntireRow.Delete()
$ws.Range("A74:A75").EntireRow.Delete()
$ws.Range("A68:A69").EntireRow.Delete()
$ws.Range("A62:A63").EntireRow.Delete()
$ws.Range("A56:A57").EntireRow.Delete()
$ws.Range("A50:A51").EntireRow.Delete()
$ws.Range("A44:A45").EntireRow.Delete()
$ws.Range("A38:A39").EntireRow.Delete()
$ws.Range("A32:A33").EntireRow.Delete()

$ws.Range("D24").Value = 8

$ws.Activate()
$selResult = $ws.Range("B174").Select()
$excel.ActiveWindow.ScrollRow = 143
